$d = $word.ActiveDocument

# Map of exact (old) paragraph text -> new (Korean) paragraph text.
# Using exact whole-paragraph text matching (rather than a document-wide
# Find/Replace) so that the identical "English" heading that lives inside
# the table-of-contents hyperlink (paragraph 1) is left untouched, while
# the "English" section heading (style P68B1DB1-Normal2) is translated.
$replacements = @{
    "English" = "영어";
    "Don’t forget to send your documents" = "잊지 않고 문서를 제출해 주시기 바랍니다";
    "If you have any questions, please contact your country manager." = "궁금하신 사항은 귀하의 국가 담당자에게 문의해 주시기 바랍니다.";
    "We look forward to seeing you there!" = "행사장에서 만나 뵙기를 기대합니다!"
}

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    # Paragraph ranges include a trailing carriage return - strip it for comparison.
    $trimmed = $text.TrimEnd([char]13, [char]7)

    if ($trimmed -eq "English") {
        # Only retarget the plain-text section heading, which uses the
        # P68B1DB1-Normal2 style; the hyperlinked "English" in the
        # language-switcher line must stay as-is.
        if ($p.Style -ne $null -and $p.Style.NameLocal -eq "P68B1DB1-Normal2") {
            $p.Range.Text = $replacements["English"]
        }
    }
    elseif ($replacements.ContainsKey($trimmed)) {
        $p.Range.Text = $replacements[$trimmed]
    }
}
